$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.343.83"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.368.22"
$ws.Range("E3").Value = "  +2.61%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.40%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0811"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.112"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").Value = "2.733.13"
$ws.Range("E15").Value = "  +2.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.82%  "

$ws.Range("D17").Value = "2.371.99"
$ws.Range("E17").Value = "  +2.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.813"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "

$ws.Range("D19").Value = "43.289.97"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("E20").Value = "  -4.61%  "

$ws.Range("D21").Value = "0.0₃0922"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("E25").Value = "  +1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("E27").Value = "  -0.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.63%  "

$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.20%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.18%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.59"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "161.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.05%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.10%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.75%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.61%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.43%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0741"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.99%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.106"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.69%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.38%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.001.47"
$ws.Range("E46").Value = "  +1.45%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0290"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.23%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.10%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.20%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.34%  "
